$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 4 — "New York -- New York": fill in the previously-empty B4:L4 cells,
# flip J4 to TRUE, and replace the rate-limit error text in O4 with "Success!"
# ---------------------------------------------------------------------------

# B4 is a date; reuse the existing date number-format so it binds to the same
# style as the other "Date Published" cells (e.g. B2/B3) instead of minting a
# new one.
$ws.Range("B4").NumberFormat = "YYYY-MM-DD"
$ws.Range("B4").Value = 44034

# C4 / D4 keep the source data's "numbers recorded as text" quirk (same as
# e.g. C8/D8 already in the sheet): force Text, write the value, then drop
# back to the default style so no new style is attached to the cell.
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "219128"
$ws.Range("C4").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "18803"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = 33790
$ws.Range("F4").Value = 5239
$ws.Range("G4").Value = 30.07
$ws.Range("H4").Value = 30.43

$ws.Range("J4").Value = $true

$ws.Range("K4").Value = 112360
$ws.Range("L4").Value = 17217

$ws.Range("O4").Value = "Success!"

# ---------------------------------------------------------------------------
# Row 41 — "Iowa": refreshed counts from the later run.
# ---------------------------------------------------------------------------
$ws.Range("C41").Value = 40146
$ws.Range("E41").Value = 3289
$ws.Range("G41").Value = 8.19

# ---------------------------------------------------------------------------
# Row 44 — "Idaho": this row now failed (timeout) instead of succeeding, so
# all the previously-populated stats are cleared back to empty text cells.
# ---------------------------------------------------------------------------
$ws.Range("B44").Value = "'"
$ws.Range("B44").Style = "Normal"

$ws.Range("C44").Value = "'"
$ws.Range("C44").Style = "Normal"

$ws.Range("D44").Value = "'"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = "'"
$ws.Range("E44").Style = "Normal"

$ws.Range("F44").Value = "'"
$ws.Range("F44").Style = "Normal"

$ws.Range("G44").Value = "'"
$ws.Range("G44").Style = "Normal"

$ws.Range("H44").Value = "'"
$ws.Range("H44").Style = "Normal"

$ws.Range("J44").Value = $false

$ws.Range("O44").Value = "An error occurred. ... TimeoutException('', None, None)"
